# Apply updated "dSF" (column F) values following a data repull / mean recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 6
    3  = -4
    4  = -2
    5  = -5
    7  = -3
    8  = -6
    11 = 1
    12 = -3
    16 = -6
    19 = 6
    20 = -1
    23 = -5
    26 = -9
    27 = -2
    28 = -5
    31 = -3
    32 = -7
    40 = -11
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
